# TradingModel - 2021/11/15 data update
# Adds a new data row (2021-11-15 / 44515, TotalProfit -1647.2) and
# shifts the "date-only" number format from the previous last row (A5)
# to the new last row (A6), matching the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (2021-11-12) was the last row before; it now uses the same
# datetime number format as the other non-final rows (s=2 in the OOXML).
$ws.Range("A5").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 6: 2021-11-15, TotalProfit -1647.2
$ws.Range("A6").Value = 44515
$ws.Range("A6").NumberFormat = "YYYY-MM-DD"
$ws.Range("B6").Value = -1647.2
